# Generate Report for Handback
# Updates the handback-status report: the "mt" (machine-translation) status
# for the second file (6ffcded5...) replaces its former "ht" (human
# translation) status, and the timestamps for that file's generate /
# handoff / handback events move forward to reflect the regenerated report.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" for 6ffcded5...md ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-28 10:17:46"
$wsOverview.Range("G4").Value = "2016-08-28 10:17:46"

# --- zh-cn sheet: Priority + Correspond Handoff/Handback DateTime ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-28 10:17:41"
$wsZhCn.Range("H4").Value = "2016-08-28 10:17:41"
$wsZhCn.Range("K3").Value = "2016-08-28 10:18:13"
$wsZhCn.Range("K4").Value = "2016-08-28 10:18:13"

# --- de-de sheet: Priority + Correspond Handoff/Handback DateTime ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-08-28 10:17:46"
$wsDeDe.Range("H4").Value = "2016-08-28 10:17:46"
$wsDeDe.Range("K3").Value = "2016-08-28 10:18:21"
$wsDeDe.Range("K4").Value = "2016-08-28 10:18:21"
